$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 379

$newData = @(
    ,@(1, "conditioned", 332, 0)
    ,@(1, "unconditioned", 332, 6)
    ,@(2, "conditioned", 332, 3)
    ,@(2, "unconditioned", 332, 3)
    ,@(3, "conditioned", 332, 0)
    ,@(3, "unconditioned", 332, 3)
    ,@(4, "conditioned", 332, 4)
    ,@(4, "unconditioned", 332, 0)
    ,@(5, "conditioned", 332, 0)
    ,@(5, "unconditioned", 332, 0)
    ,@(6, "conditioned", 332, 0)
    ,@(6, "unconditioned", 332, 1)
    ,@(7, "conditioned", 332, 0)
    ,@(7, "unconditioned", 332, 0)
    ,@(8, "conditioned", 332, 0)
    ,@(8, "unconditioned", 332, 1)
    ,@(9, "conditioned", 332, 0)
    ,@(9, "unconditioned", 332, 0)
    ,@(10, "conditioned", 332, 1)
    ,@(10, "unconditioned", 332, 1)
    ,@(11, "conditioned", 332, 0)
    ,@(11, "unconditioned", 332, 0)
    ,@(12, "conditioned", 332, 2)
    ,@(12, "unconditioned", 332, 1)
    ,@(13, "conditioned", 332, 1)
    ,@(13, "unconditioned", 332, 1)
    ,@(14, "conditioned", 332, 2)
    ,@(14, "unconditioned", 332, 0)
    ,@(15, "unconditioned", 332, $null)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($null -ne $row[3]) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
}

$ws.Range("D407").Select() | Out-Null

